$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new song entries (updated 2018-11-16), following the pattern of
# the existing rows. Copy formatting from the last existing data row (45)
# so the new "#" column (A) keeps its bold/centered/bordered style.

# Row 46: "Sad_" by XXXTENTACION
$ws.Cells.Item(45, 1).Copy($ws.Cells.Item(46, 1))
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = "Sad_"
$ws.Cells.Item(46, 3).Value = "XXXTENTACION"
$ws.Cells.Item(46, 4).Value = "16-11-2018"

# Row 47: "No Brainer" by DJ Khaled, Justin Bieber, C
$ws.Cells.Item(45, 1).Copy($ws.Cells.Item(47, 1))
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = "No Brainer"
$ws.Cells.Item(47, 3).Value = "DJ Khaled, Justin Bieber, C"
$ws.Cells.Item(47, 4).Value = "16-11-2018"
